$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name / title
$ws.Name = "Through 2021-11-14"

# Update the "November (through 11-13)" label to "November (through 11-14)"
$ws.Range("A12").Value = "November (through 11-14)"

# Update November row (row 12) values
$ws.Range("B12").Value = 17
$ws.Range("C12").Value = 35
$ws.Range("D12").Value = 62
$ws.Range("E12").Value = 29
$ws.Range("F12").Value = 23
$ws.Range("G12").Value = 87
$ws.Range("H12").Value = 97

# Update Total row (row 13) values
$ws.Range("B13").Value = 275
$ws.Range("C13").Value = 521
$ws.Range("D13").Value = 772
$ws.Range("E13").Value = 644
$ws.Range("F13").Value = 505
$ws.Range("G13").Value = 1144
$ws.Range("H13").Value = 1539
